# Update the "Metadata" sheet of the ConceptMap spreadsheet export:
#   - URL value gets a path tweak (.../fhir/fr/medication/... -> .../ig/fhir/medication/...)
#   - Date value bumped to the new publication timestamp
#   - Jurisdiction value filled in with "FRANCE" (was blank)
#   - Target value gets the same path tweak as the URL above

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ConceptMap/PN13-FHIR-prescmed-patient-id-seul-conceptmap"
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$ws.Range("B11").Value = "FRANCE"
$ws.Range("B16").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-inpatient-medicationrequest"
